# Update the "Plans" worksheet (residential plan pricing table).
#
# The existing 5 rows of the "Residencial" / "Sin_TotalPlay_TV" plan had
# their "Megas" values halved, and a 6th row was added for a new
# "Residencial" / "Con_TotalPlay_TV" (3P BRM) plan at 1000 Megas.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plans")
$ws.Activate()

# Halve the Megas values for the existing Sin_TotalPlay_TV residential plans.
$ws.Range("D5").Value = 20
$ws.Range("D6").Value = 50
$ws.Range("D7").Value = 100
$ws.Range("D8").Value = 200
$ws.Range("D9").Value = 500

# New plan row: Residencial / Con_TotalPlay_TV (3P BRM) at 1000 Megas.
$ws.Range("B10").Value = "Residencial"
$ws.Range("C10").Value = "Con_TotalPlay_TV"
$ws.Range("D10").Value = 1000

# Restore the author's cursor position after the edit.
$ws.Range("E14").Select()
